# Final updates after clean run
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Program")

# Update cell values with the new text content
$ws.Range("A5").Value = "Team4kl8Team4Team4"
$ws.Range("B5").Value = "javajavkli09Team4"
$ws.Range("A6").Value = "JavakjhgflTeam4"
$ws.Range("B7").Value = "XkjmnhjkmnhjTeam4"

# Make the Program sheet active and move the selection to B7
$ws.Activate()
$ws.Range("B7").Select()
